$d = $word.ActiveDocument

# 1. Merge the "Improvements for V0.2 of the [Rastaban] PCB" title runs
#    (removes the spell-check proofErr markers around "Rastaban" and
#    collapses the three runs into a single run) via a find/replace of
#    the whole paragraph text with itself.
$titleRange = $d.Paragraphs(1).Range
$null = $titleRange.Find.Execute("Improvements for V0.2 of the Rastaban PCB", $true, $false, $false, $false, $false, $true, 1, $false, "Improvements for V0.2 of the Rastaban PCB", 2)

# 2. Add the new bullet point text to the previously-empty list paragraph.
#    The text is typed as several distinct runs (mirrors separate typing/
#    editing actions), so each chunk is inserted as its own paragraph and
#    then stitched back together by deleting the paragraph mark between
#    them -- this keeps the runs from being silently coalesced into one.
$parts = @(
    "For the 2209 ",
    "we could use",
    " ",
    "UART port",
    ", so ",
    "we",
    " can control clock, microstepping etc, set parameters, delays",
    ", coolstepping, ",
    "etc."
)

$targetIndex = 4
$p = $d.Paragraphs($targetIndex)
$insertPoint = $d.Range($p.Range.Start, $p.Range.Start)
$insertPoint.InsertAfter($parts[0])

for ($i = 1; $i -lt $parts.Length; $i++) {
    $p = $d.Paragraphs($targetIndex)
    $endOfPara = $d.Range($p.Range.End - 1, $p.Range.End - 1)
    $endOfPara.InsertParagraphAfter()

    $newPara = $d.Paragraphs($targetIndex + 1)
    $newPara.Range.InsertAfter($parts[$i])

    $mergedPara = $d.Paragraphs($targetIndex)
    $paraMark = $d.Range($mergedPara.Range.End - 1, $mergedPara.Range.End)
    $paraMark.Delete()
}

Write-Output "done"
